$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new "2026-1" worksheet as the very first tab, cloned from the current
# first sheet ("2025-2") so that it inherits the same column widths / cell
# styles, then overwrite its contents with the new season's data.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(1)
$template.Copy($template) | Out-Null
$ws = $wb.Worksheets.Item(1)
$ws.Name = "2026-1"

# New data rows for the 2026-1 season.
$rows = @(
    @("2026-1", "CHATA GUANAY", "Chata", "GO/50", "GO/50-126"),
    @("2026-1", "CHATA TASA", "Chata", "GO/51", "GO/51-126"),
    @("2026-1", "REM LOBOS", "Remolcador", "A.S/0055", "A.S/0055-126"),
    @("2026-1", "EQUIPOS PRUBA DE ESTABILIDAD", "Embarcación Pesquera", "GP/94", "GP/94-126"),
    @("2026-1", "EP TASA 414", "Embarcación Pesquera", "GP/99", "GP/99-126"),
    @("2026-1", "EP TASA 411", "Embarcación Pesquera", "GP/97", "GP/97-126")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Rows 8-19 keep only the (pre-existing) style on column C, with no values -
# clear out the data that was copied from the template sheet.
$ws.Range("A8:E19").ClearContents()

# The previously-active "Modulo" tab loses its selected/tabbed state, but its
# recorded cursor position moved too (A5 -> D10) - replicate that, then
# re-select the new sheet last so it ends up the active tab on save.
$modulo = $wb.Worksheets.Item("Modulo")
$modulo.Range("D10").Select() | Out-Null

# Restore the selection/active-cell state recorded for the new tab.
$ws.Select() | Out-Null
$ws.Range("E11").Select() | Out-Null
